# Guestbook entries: update row 2 with the real submitted values (the
# placeholder test row "asdas" / "181485892" / ... gets replaced) and
# append a new guestbook entry in row 3.
#
# Columns: A=name, B=rut, C=age, D=phone, E=message, F=timestamp
#
# Some of these values look numeric ("21", "986389894", ...) but the
# source data is plain text, so a leading apostrophe is used to force
# Excel to store them as text instead of silently converting them to
# numbers. The Style reset afterwards clears the quote-prefix
# formatting flag that Excel applies to the cell when you do that,
# leaving the cell on the sheet's default (unstyled) look - matching
# how the original rows were written out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the placeholder entry with the real one ---
$ws.Range("A2").Value = "Andrea Castillo"
$ws.Range("B2").Value = "20985370-1"

$ws.Range("C2").Value = "'21"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "'986389894"
$ws.Range("D2").Style = "Normal"

# E2 (message) is already blank/empty text - leave it as-is.

$ws.Range("F2").Value = "2023-11-24 10:58:03"

# --- Row 3: new guestbook entry ---
$ws.Range("A3").Value = "Sofia Martinez"
$ws.Range("B3").Value = "11161499-7"

$ws.Range("C3").Value = "'21"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "'934393434"
$ws.Range("D3").Style = "Normal"

# No message was submitted - keep it an empty text cell (matches E2).
$ws.Range("E3").Value = "'"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = "2023-11-24 10:58:31"
